# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the combined "全部类型" sheet, matching the newly scraped totals.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - rows 3-19
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 38
$wsExpo.Range("F4").Value = 1417
$wsExpo.Range("F5").Value = 324
$wsExpo.Range("F6").Value = 1039
$wsExpo.Range("F10").Value = 298
$wsExpo.Range("F11").Value = 1043
$wsExpo.Range("F12").Value = 721
$wsExpo.Range("F13").Value = 12093
$wsExpo.Range("F14").Value = 12559
$wsExpo.Range("F19").Value = 76

# Sheet "全部类型" (all types) - same events, offset by one row
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 38
$wsAll.Range("F5").Value = 1417
$wsAll.Range("F6").Value = 324
$wsAll.Range("F7").Value = 1039
$wsAll.Range("F11").Value = 298
$wsAll.Range("F12").Value = 1043
$wsAll.Range("F13").Value = 721
$wsAll.Range("F14").Value = 12093
$wsAll.Range("F15").Value = 12559
$wsAll.Range("F20").Value = 76
